$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 920, shifting existing rows 920:990 down to 922:992
$ws.Range("A920:R921").Insert(-4121)  # xlShiftDown = -4121

# Populate the two new rows (920 and 921) with the new records.
# Row 920
$ws.Cells.Item(920, 1).Value = 11
$ws.Cells.Item(920, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(920, 3).Value = "Bíobío"
$ws.Cells.Item(920, 4).Value = 45265
$ws.Cells.Item(920, 5).Value = 8
$ws.Cells.Item(920, 6).Value = 100112004
$ws.Cells.Item(920, 7).Value = "Cebolla"
$ws.Cells.Item(920, 8).Value = "Sin especificar"
$ws.Cells.Item(920, 9).Value = "1a (cosecha)"
$ws.Cells.Item(920, 10).Value = 500
$ws.Cells.Item(920, 11).Value = 13000
$ws.Cells.Item(920, 12).Value = 14000
$ws.Cells.Item(920, 13).Value = 13600
$ws.Cells.Item(920, 14).Value = "`$/malla 16 kilos"
$ws.Cells.Item(920, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(920, 16).Value = 850
$ws.Cells.Item(920, 17).Value = 16
$ws.Cells.Item(920, 18).Value = "Hortaliza"

# Row 921
$ws.Cells.Item(921, 1).Value = 11
$ws.Cells.Item(921, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(921, 3).Value = "Bíobío"
$ws.Cells.Item(921, 4).Value = 45265
$ws.Cells.Item(921, 5).Value = 8
$ws.Cells.Item(921, 6).Value = 100112004
$ws.Cells.Item(921, 7).Value = "Cebolla"
$ws.Cells.Item(921, 8).Value = "Sin especificar"
$ws.Cells.Item(921, 9).Value = "2a (cosecha)"
$ws.Cells.Item(921, 10).Value = 200
$ws.Cells.Item(921, 11).Value = 12000
$ws.Cells.Item(921, 12).Value = 12000
$ws.Cells.Item(921, 13).Value = 12000
$ws.Cells.Item(921, 14).Value = "`$/malla 16 kilos"
$ws.Cells.Item(921, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(921, 16).Value = 750
$ws.Cells.Item(921, 17).Value = 16
$ws.Cells.Item(921, 18).Value = "Hortaliza"
